$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (contrast = NN effect)
$ws.Range("B2").Value = 0.0898550975333272
$ws.Range("C2").Value = 0.0316886475803865
$ws.Range("E2").Value = 0.0277464895569875
$ws.Range("F2").Value = 0.151963705509667
$ws.Range("G2").Value = 2.83556113606258
$ws.Range("H2").Value = 0.0045745258726248

# Row 3 (contrast = Used_CAM_ONLY effect)
$ws.Range("B3").Value = -0.107889000538331
$ws.Range("C3").Value = 0.0720581427078311
$ws.Range("E3").Value = -0.249120365038528
$ws.Range("F3").Value = 0.0333423639618653
$ws.Range("G3").Value = -1.49724925572646
$ws.Range("H3").Value = 0.134328413240054

# Row 4 (contrast = Used_conv_and_CAM effect)
$ws.Range("B4").Value = -0.0322218945184243
$ws.Range("C4").Value = 0.0317183223604764
$ws.Range("E4").Value = -0.0943886639949896
$ws.Range("F4").Value = 0.029944874958141
$ws.Range("G4").Value = -1.01587638060503
$ws.Range("H4").Value = 0.309688257683836

# Row 5 (contrast = Used_conv_ONLY effect)
$ws.Range("B5").Value = 0.0502557975234283
$ws.Range("C5").Value = 0.0268858181835791
$ws.Range("E5").Value = -0.00243943781127885
$ws.Range("F5").Value = 0.102951032858136
$ws.Range("G5").Value = 1.86923072901396
$ws.Range("H5").Value = 0.0615907207859216
